# Updated cryptos list on Sat Apr  1 10:58:09 UTC 2023 with GitHub Actions
#
# For each row in the "cryptos" table, refresh the Price (column D) and
# Volume(1h) (column E) values with the latest scraped figures.
#
# Column D values are plain text in the workbook (e.g. "28.482.08",
# "1.110", "0.07664") even though many of them look numeric. Excel's COM
# layer will happily auto-coerce a single-decimal-point string like
# "41.86" into a real number (and normalise away significant trailing
# zeros, e.g. "1.000" -> 1), so before writing any such value we force
# the destination cell to Text format ("@") first. Values that contain
# two '.' separators (thousand/decimal grouping, e.g. "28.482.08") can
# never be parsed as a number by Excel, so that extra step is skipped
# for those to avoid needlessly touching their style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  Price = "28.482.08";   Volume = "  +2.00%  " },
    @{ Row = 3;  Price = "1.828.69";    Volume = "  +1.83%  " },
    @{ Row = 4;  Price = $null;         Volume = "  -0.11%  " },
    @{ Row = 5;  Price = "315.37";      Volume = "  -0.52%  " },
    @{ Row = 6;  Price = $null;         Volume = "  +0.01%  " },
    @{ Row = 7;  Price = "0.5077";      Volume = "  -4.50%  " },
    @{ Row = 8;  Price = "0.3906";      Volume = "  +0.64%  " },
    @{ Row = 9;  Price = "0.07664";     Volume = "  +3.03%  " },
    @{ Row = 10; Price = "41.86";       Volume = "  +1.04%  " },
    @{ Row = 11; Price = "1.110";       Volume = "  +2.16%  " },
    @{ Row = 12; Price = "21.12";       Volume = "  +3.80%  " },
    @{ Row = 13; Price = "6.297";       Volume = "  +2.03%  " },
    @{ Row = 14; Price = "7.593";       Volume = "  +2.19%  " },
    @{ Row = 15; Price = "1.000";       Volume = "  -0.13%  " },
    @{ Row = 16; Price = "1.822.60";    Volume = "  +1.43%  " },
    @{ Row = 17; Price = "93.26";       Volume = "  +5.56%  " },
    @{ Row = 18; Price = "0.00001084";  Volume = "  +2.35%  " },
    @{ Row = 19; Price = "0.06679";     Volume = "  +2.15%  " },
    @{ Row = 20; Price = "17.71";       Volume = "  +2.80%  " },
    @{ Row = 21; Price = $null;         Volume = "  +0.04%  " },
    @{ Row = 22; Price = "6.149";       Volume = "  +3.42%  " },
    @{ Row = 23; Price = "28.505.98";   Volume = "  +1.93%  " },
    @{ Row = 24; Price = "11.14";       Volume = "  +0.43%  " },
    @{ Row = 25; Price = $null;         Volume = "  +7.88%  " },
    @{ Row = 26; Price = "156.94";      Volume = "  -0.27%  " },
    @{ Row = 27; Price = "20.62";       Volume = "  +2.49%  " },
    @{ Row = 28; Price = "2.034.22";    Volume = "  +1.66%  " },
    @{ Row = 29; Price = "2.400";       Volume = "  +4.74%  " },
    @{ Row = 30; Price = "125.37";      Volume = "  +2.93%  " },
    @{ Row = 31; Price = "1.130";       Volume = "  +2.99%  " },
    @{ Row = 32; Price = "0.1084";      Volume = "  -0.39%  " },
    @{ Row = 33; Price = "5.687";       Volume = "  +3.59%  " },
    @{ Row = 34; Price = "3.661";       Volume = "  -0.17%  " },
    @{ Row = 35; Price = "0.07049";     Volume = "  +0.27%  " },
    @{ Row = 36; Price = "0.2234";      Volume = "  +1.46%  " },
    @{ Row = 37; Price = "8.976";       Volume = "  +7.08%  " },
    @{ Row = 38; Price = "0.02325";     Volume = "  +2.18%  " },
    @{ Row = 39; Price = "5.145";       Volume = "  +1.36%  " },
    @{ Row = 40; Price = "0.6268";      Volume = "  +2.64%  " },
    @{ Row = 41; Price = "11.24";       Volume = "  +0.30%  " },
    @{ Row = 42; Price = "1.181";       Volume = "  -0.32%  " },
    @{ Row = 43; Price = $null;         Volume = "  -0.01%  " },
    @{ Row = 44; Price = "1.397";       Volume = "  -1.52%  " },
    @{ Row = 45; Price = "13.40";       Volume = "  +0.25%  " },
    @{ Row = 46; Price = "0.5907";      Volume = "  +3.67%  " },
    @{ Row = 47; Price = "3.716";       Volume = "  +1.11%  " },
    @{ Row = 48; Price = "124.56";      Volume = "  -0.19%  " },
    @{ Row = 49; Price = "1.983";       Volume = "  +3.60%  " },
    @{ Row = 50; Price = "1.195";       Volume = "  +1.59%  " },
    @{ Row = 51; Price = "0.06919";     Volume = "  +1.64%  " }
)

foreach ($item in $rows) {
    $r = $item.Row

    if ($null -ne $item.Price) {
        $priceCell = $ws.Range("D$r")
        # Only single '.' numeric-looking strings risk being silently
        # reinterpreted as a Number by Excel; two-dot strings like
        # "28.482.08" are never valid numbers, so leave their (absent)
        # number format untouched to avoid unrelated style churn.
        if (([regex]::Matches($item.Price, "\.")).Count -lt 2) {
            $priceCell.NumberFormat = "@"
        }
        $priceCell.Value = $item.Price
    }

    $ws.Range("E$r").Value = $item.Volume
}
